# data_summary.xlsx — append new rows of video-derived stats and correct
# the two most-recent hourly rows (hour 29 -> 45, hour 30 -> 46, etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "hour" column (B) holds text-looking numbers (e.g. "29"), not real
# numbers, in the original sheet. A leading apostrophe forces Excel to
# keep the entry as text instead of auto-converting it, and resetting the
# style back to "Normal" afterwards clears the quote-prefix formatting
# flag that the apostrophe entry leaves behind.

# Row 2: hour 29 -> 45, #people 37 -> 26, #males 10 -> 6 (#females, age unchanged)
$ws.Range("B2").Value = "'45"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 26
$ws.Range("D2").Value = 6

# Row 3: hour 30 -> 46, #males 15 -> 14, #females 3 -> 2, age 15-25 -> 6-14
$ws.Range("B3").Value = "'46"
$ws.Range("B3").Style = "Normal"
$ws.Range("D3").Value = 14
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = "6-14"

# Row 4 (new): 29/05/2022, hour 47, 57 people, 18 males, 6 females, 15-25
$ws.Range("A4").Value = "29/05/2022"
$ws.Range("B4").Value = "'47"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 57
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = "15-25"

# Row 5 (new): 29/05/2022, hour 48, 44 people, 12 males, 2 females, 15-25
$ws.Range("A5").Value = "29/05/2022"
$ws.Range("B5").Value = "'48"
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").Value = 44
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = "15-25"

# Row 6 (new): 29/05/2022, hour 49, 78 people, 7 males, 4 females, 6-14
$ws.Range("A6").Value = "29/05/2022"
$ws.Range("B6").Value = "'49"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 78
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = "6-14"
